# Auto-generated Excel COM-interop script
# Applies numeric value updates to the Leve profit-tracking sheets
# (currentAveragePrice* / LevePrice* / LeveProfit* columns H-N)
# to match refreshed Universalis market data.

$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")

# Row 2
$ws.Range("H2").Value = 368.1
$ws.Range("I2").Value = 355.16666
$ws.Range("K2").Value = 355.16666
$ws.Range("M2").Value = -242.16666

# Row 33
$ws.Range("H33").Value = 435.46667
$ws.Range("I33").Value = 435.46667
$ws.Range("K33").Value = 435.46667
$ws.Range("M33").Value = -206.46667

# Row 40
$ws.Range("H40").Value = 1676.1428
$ws.Range("I40").Value = 1418.1818
$ws.Range("K40").Value = 1418.1818
$ws.Range("M40").Value = -1243.1818

# Row 42
$ws.Range("H42").Value = 363.26666
$ws.Range("I42").Value = 278.6
$ws.Range("J42").Value = 405.6
$ws.Range("K42").Value = 835.8000000000001
$ws.Range("L42").Value = 1216.8
$ws.Range("M42").Value = -605.8000000000001
$ws.Range("N42").Value = -1676.8

# Row 64
$ws.Range("H64").Value = 3475.6
$ws.Range("I64").Value = 3444.4614
$ws.Range("J64").Value = 3499.4119
$ws.Range("K64").Value = 3444.4614
$ws.Range("L64").Value = 3499.4119
$ws.Range("M64").Value = -3196.4614
$ws.Range("N64").Value = -3995.4119

# Row 67
$ws.Range("H67").Value = 3475.6
$ws.Range("I67").Value = 3444.4614
$ws.Range("J67").Value = 3499.4119
$ws.Range("K67").Value = 3444.4614
$ws.Range("L67").Value = 3499.4119
$ws.Range("M67").Value = -2586.4614
$ws.Range("N67").Value = -5215.4119

# Row 132
$ws.Range("H132").Value = 2036.7959
$ws.Range("I132").Value = 1798.8
$ws.Range("J132").Value = 3094.5557
$ws.Range("K132").Value = 5396.4
$ws.Range("L132").Value = 9283.667099999999
$ws.Range("M132").Value = -2866.4
$ws.Range("N132").Value = -14343.6671

# Row 137
$ws.Range("H137").Value = 1357.9117
$ws.Range("I137").Value = 1295.8064
$ws.Range("J137").Value = 1999.6666
$ws.Range("K137").Value = 3887.4192
$ws.Range("L137").Value = 5998.9998
$ws.Range("M137").Value = -1337.4192
$ws.Range("N137").Value = -11098.9998


# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")

# Row 12
$ws.Range("H12").Value = 2500
$ws.Range("I12").Value = 2500
$ws.Range("K12").Value = 2500
$ws.Range("M12").Value = -2327

# Row 60
$ws.Range("H60").Value = 50000
$ws.Range("I60").Value = 50000
$ws.Range("K60").Value = 50000
$ws.Range("M60").Value = -49267

# Row 61
$ws.Range("H61").Value = 1423.186
$ws.Range("I61").Value = 1208.4857
$ws.Range("K61").Value = 1208.4857
$ws.Range("M61").Value = -996.4857

# Row 74
$ws.Range("H74").Value = 1306.8462
$ws.Range("I74").Value = 916.5
$ws.Range("K74").Value = 916.5
$ws.Range("M74").Value = -42.5

# Row 77
$ws.Range("H77").Value = 1306.8462
$ws.Range("I77").Value = 916.5
$ws.Range("K77").Value = 4582.5
$ws.Range("M77").Value = -214.5

# Row 136
$ws.Range("H136").Value = 1423.186
$ws.Range("I136").Value = 1208.4857
$ws.Range("K136").Value = 3625.4571
$ws.Range("M136").Value = -1075.4571


# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")

# Row 32
$ws.Range("H32").Value = 52500
$ws.Range("J32").Value = 55000
$ws.Range("L32").Value = 55000
$ws.Range("N32").Value = -55768


# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")

# Row 132
$ws.Range("H132").Value = 308710.5
$ws.Range("I132").Value = 376656.94
$ws.Range("J132").Value = 2951.5
$ws.Range("K132").Value = 1129970.82
$ws.Range("L132").Value = 8854.5
$ws.Range("M132").Value = -1127440.82
$ws.Range("N132").Value = -13914.5

# Row 134
$ws.Range("H134").Value = 1861.1041
$ws.Range("I134").Value = 1284.8889
$ws.Range("J134").Value = 3589.75
$ws.Range("K134").Value = 3854.6667
$ws.Range("L134").Value = 10769.25
$ws.Range("M134").Value = -1319.6667
$ws.Range("N134").Value = -15839.25


# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")

# Row 21
$ws.Range("H21").Value = 3032944
$ws.Range("I21").Value = 683
$ws.Range("J21").Value = 4170041.8
$ws.Range("K21").Value = 2049
$ws.Range("L21").Value = 12510125.4
$ws.Range("M21").Value = -1876
$ws.Range("N21").Value = -12510471.4

# Row 63
$ws.Range("H63").Value = 526881
$ws.Range("I63").Value = 526881
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 1580643
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -1579894
$ws.Range("N63").ClearContents()

# Row 66
$ws.Range("H66").Value = 526881
$ws.Range("I66").Value = 526881
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 4741929
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -4738185
$ws.Range("N66").ClearContents()

# Row 70
$ws.Range("H70").Value = 8214
$ws.Range("I70").Value = 10742.4
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 32227.2
$ws.Range("L70").Value = 12000
$ws.Range("M70").Value = -31912.2
$ws.Range("N70").Value = -12630

# Row 73
$ws.Range("H73").Value = 8214
$ws.Range("I73").Value = 10742.4
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 32227.2
$ws.Range("L73").Value = 12000
$ws.Range("M73").Value = -31135.2
$ws.Range("N73").Value = -14184

# Row 75
$ws.Range("H75").Value = 9166.666999999999
$ws.Range("I75").Value = 1000
$ws.Range("J75").Value = 10800
$ws.Range("K75").Value = 3000
$ws.Range("L75").Value = 32400
$ws.Range("M75").Value = -2002
$ws.Range("N75").Value = -34396

# Row 78
$ws.Range("H78").Value = 9166.666999999999
$ws.Range("I78").Value = 1000
$ws.Range("J78").Value = 10800
$ws.Range("K78").Value = 9000
$ws.Range("L78").Value = 97200
$ws.Range("M78").Value = -4008
$ws.Range("N78").Value = -107184

# Row 87
$ws.Range("H87").Value = 8648.387000000001
$ws.Range("I87").Value = 1525
$ws.Range("J87").Value = 9703.704
$ws.Range("K87").Value = 4575
$ws.Range("L87").Value = 29111.112
$ws.Range("M87").Value = -3327
$ws.Range("N87").Value = -31607.112

# Row 90
$ws.Range("H90").Value = 8648.387000000001
$ws.Range("I90").Value = 1525
$ws.Range("J90").Value = 9703.704
$ws.Range("K90").Value = 13725
$ws.Range("L90").Value = 87333.336
$ws.Range("M90").Value = -7485
$ws.Range("N90").Value = -99813.336

# Row 92
$ws.Range("H92").Value = 470
$ws.Range("I92").Value = 470
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 1410
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -162
$ws.Range("N92").ClearContents()

# Row 109
$ws.Range("H109").Value = 3345.8333
$ws.Range("I109").Value = 750
$ws.Range("J109").Value = 3865
$ws.Range("K109").Value = 2250
$ws.Range("L109").Value = 11595
$ws.Range("M109").Value = -1210
$ws.Range("N109").Value = -13675

# Row 114
$ws.Range("H114").Value = 953.45
$ws.Range("I114").Value = 333.7
$ws.Range("J114").Value = 1573.2
$ws.Range("K114").Value = 1001.1
$ws.Range("L114").Value = 4719.6
$ws.Range("M114").Value = 2252.9
$ws.Range("N114").Value = -11227.6

# Row 117
$ws.Range("H117").Value = 672.3
$ws.Range("I117").Value = 477.8
$ws.Range("J117").Value = 737.13336
$ws.Range("K117").Value = 1433.4
$ws.Range("L117").Value = 2211.40008
$ws.Range("M117").Value = 2008.6
$ws.Range("N117").Value = -9095.400079999999

# Row 131
$ws.Range("H131").Value = 1025.96
$ws.Range("I131").Value = 426
$ws.Range("J131").Value = 1057.5369
$ws.Range("K131").Value = 1278
$ws.Range("L131").Value = 3172.6107
$ws.Range("M131").Value = 3762
$ws.Range("N131").Value = -13252.6107

# Row 132
$ws.Range("H132").Value = 1621.8636
$ws.Range("I132").Value = 1249.7142
$ws.Range("J132").Value = 1795.5333
$ws.Range("K132").Value = 11247.4278
$ws.Range("L132").Value = 16159.7997
$ws.Range("M132").Value = -8717.427799999999
$ws.Range("N132").Value = -21219.7997

# Row 139
$ws.Range("H139").Value = 2368.6785
$ws.Range("I139").Value = 1891.25
$ws.Range("J139").Value = 3005.25
$ws.Range("K139").Value = 5673.75
$ws.Range("L139").Value = 9015.75
$ws.Range("M139").Value = -533.75
$ws.Range("N139").Value = -19295.75


# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")

# Row 12
$ws.Range("H12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("N12").ClearContents()

# Row 20
$ws.Range("H20").Value = 42501.5
$ws.Range("J20").Value = 40002
$ws.Range("L20").Value = 40002
$ws.Range("N20").Value = -40492

# Row 52
$ws.Range("H52").Value = 50000
$ws.Range("I52").Value = 50000
$ws.Range("J52").Value = 50000
$ws.Range("K52").Value = 50000
$ws.Range("L52").Value = 50000
$ws.Range("M52").Value = -49741
$ws.Range("N52").Value = -50518

# Row 58
$ws.Range("H58").Value = 32000
$ws.Range("I58").Value = 3000
$ws.Range("J58").Value = 41666.668
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 41666.668
$ws.Range("M58").Value = -2723
$ws.Range("N58").Value = -42220.668

# Row 107
$ws.Range("H107").Value = 2982
$ws.Range("I107").Value = 3709.3333
$ws.Range("K107").Value = 3709.3333
$ws.Range("M107").Value = -1789.3333


# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")

# Row 44
$ws.Range("H44").Value = 4117.25
$ws.Range("J44").Value = 4117.25
$ws.Range("L44").Value = 4117.25
$ws.Range("N44").Value = -5029.25

# Row 132
$ws.Range("H132").Value = 5418.6875
$ws.Range("I132").Value = 5361.615
$ws.Range("K132").Value = 16084.845
$ws.Range("M132").Value = -13554.845


# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")

# Row 132
$ws.Range("H132").Value = 1953.6
$ws.Range("I132").Value = 1323.4762
$ws.Range("J132").Value = 2898.7856
$ws.Range("K132").Value = 3970.4286
$ws.Range("L132").Value = 8696.356800000001
$ws.Range("M132").Value = -1440.4286
$ws.Range("N132").Value = -13756.3568

# Row 136
$ws.Range("H136").Value = 2151.4211
$ws.Range("I136").Value = 1916.5
$ws.Range("J136").Value = 2554.1428
$ws.Range("K136").Value = 5749.5
$ws.Range("L136").Value = 7662.428400000001
$ws.Range("M136").Value = -3199.5
$ws.Range("N136").Value = -12762.4284

